$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "#merkuur,"
$ws.Range("C2").Value = "Merkuur,"
$ws.Range("D2").Value = ""

$ws.Range("B3").Value = "#urbanus"
$ws.Range("C3").Value = "Urbanus"
$ws.Range("D3").Value = ""

$ws.Range("B4").Value = "#merkuur"
$ws.Range("C4").Value = "Merkuur"
$ws.Range("D4").Value = ""

$ws.Range("B5").Value = "#klara"
$ws.Range("C5").Value = "Klara"
$ws.Range("D5").Value = ""

$ws.Range("B6").Value = "#isabella,"

$ws.Range("B7").Value = "#isabella"
$ws.Range("C7").Value = "Isabella"
